$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: move value from C2 to B2
$ws.Range("C2").ClearContents()
$ws.Range("B2").Value = 2.5

# Row 5: combine B5+C5 into B5, clear C5
$ws.Range("B5").Value = 5.5
$ws.Range("C5").ClearContents()

# Row 31: update B31 and C31
$ws.Range("B31").Value = 2.75
$ws.Range("C31").Value = 1

# Row 38: combine B38+C38 into B38, clear C38
$ws.Range("B38").Value = 3
$ws.Range("C38").ClearContents()

# Row 54: combine B54+C54 into B54, clear C54
$ws.Range("B54").Value = 3.25
$ws.Range("C54").ClearContents()

# Update sheet view: remove topLeftCell scroll, change selection to B5
$ws.Range("B5").Select()
